$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 45: The House Always Wins / Blinding Potion
$ws.Range("H45").Value = 4755.6665
$ws.Range("I45").Value = 4506.8
$ws.Range("K45").Value = 13520.4
$ws.Range("M45").Value = -13328.4
# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 2523.6462
$ws.Range("I138").Value = 1565.5927
$ws.Range("J138").Value = 3204.3684
$ws.Range("K138").Value = 4696.7781
$ws.Range("L138").Value = 9613.1052
$ws.Range("M138").Value = 443.2219000000005
$ws.Range("N138").Value = -19893.1052

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 7917.958
$ws.Range("I32").Value = 7144.0234
$ws.Range("J32").Value = 15313.333
$ws.Range("K32").Value = 7144.0234
$ws.Range("L32").Value = 15313.333
$ws.Range("M32").Value = -6857.0234
$ws.Range("N32").Value = -15887.333
# Row 45: Hollow Hallmarks / Mythril Ingot
$ws.Range("H45").Value = 2278.15
$ws.Range("I45").Value = 2073.4707
$ws.Range("J45").Value = 3438
$ws.Range("K45").Value = 2073.4707
$ws.Range("L45").Value = 3438
$ws.Range("M45").Value = -1696.4707
$ws.Range("N45").Value = -4192
# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 1779.8718
$ws.Range("I61").Value = 1148.963
$ws.Range("J61").Value = 3199.4167
$ws.Range("K61").Value = 1148.963
$ws.Range("L61").Value = 3199.4167
$ws.Range("M61").Value = -936.963
$ws.Range("N61").Value = -3623.4167
# Row 109: A Head of Demand / Deepgold Helm of Fending
$ws.Range("H109").Value = 45250
$ws.Range("J109").Value = 45250
$ws.Range("L109").Value = 45250
$ws.Range("N109").Value = -48024
# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 16131251
$ws.Range("I132").Value = 27779570
$ws.Range("J132").Value = 2808.2307
$ws.Range("K132").Value = 83338710
$ws.Range("L132").Value = 8424.6921
$ws.Range("M132").Value = -83336180
$ws.Range("N132").Value = -13484.6921
# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 1779.8718
$ws.Range("I136").Value = 1148.963
$ws.Range("J136").Value = 3199.4167
$ws.Range("K136").Value = 3446.889
$ws.Range("L136").Value = 9598.250100000001
$ws.Range("M136").Value = -896.8890000000001
$ws.Range("N136").Value = -14698.2501

$ws = $wb.Worksheets.Item("BSM")
# Row 6: The Unkindest Cut / Bronze Saw
$ws.Range("H6").Value = 19912.445
$ws.Range("J6").Value = 19912.445
$ws.Range("L6").Value = 19912.445
$ws.Range("N6").Value = -20138.445
# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 2850.4443
$ws.Range("I134").Value = 2410.9614
$ws.Range("J134").Value = 3993.1
$ws.Range("K134").Value = 7232.8842
$ws.Range("L134").Value = 11979.3
$ws.Range("M134").Value = -4697.8842
$ws.Range("N134").Value = -17049.3

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 3748722.2
$ws.Range("I31").Value = 1504.4
$ws.Range("J31").Value = 7581104.5
$ws.Range("K31").Value = 1504.4
$ws.Range("L31").Value = 7581104.5
$ws.Range("M31").Value = -1209.4
$ws.Range("N31").Value = -7581694.5
# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 3748722.2
$ws.Range("I34").Value = 1504.4
$ws.Range("J34").Value = 7581104.5
$ws.Range("K34").Value = 1504.4
$ws.Range("L34").Value = 7581104.5
$ws.Range("M34").Value = -1302.4
$ws.Range("N34").Value = -7581508.5
# Row 96: Composition / Larch Composite Bow
$ws.Range("H96").Value = 51952.57
$ws.Range("J96").Value = 51952.57
$ws.Range("L96").Value = 51952.57
$ws.Range("N96").Value = -57444.57
# Row 120: Kindling the Flame / Lignum Vitae Ring
$ws.Range("H120").Value = 30619.666
$ws.Range("J120").Value = 30619.666
$ws.Range("L120").Value = 30619.666
$ws.Range("N120").Value = -37877.666
# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 313280.94
$ws.Range("I132").Value = 1312.3715
$ws.Range("J132").Value = 1405170.9
$ws.Range("K132").Value = 3937.1145
$ws.Range("L132").Value = 4215512.699999999
$ws.Range("M132").Value = -1407.1145
$ws.Range("N132").Value = -4220572.699999999
# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 760303.7
$ws.Range("I134").Value = 441213.25
$ws.Range("J134").Value = 2802482.5
$ws.Range("K134").Value = 1323639.75
$ws.Range("L134").Value = 8407447.5
$ws.Range("M134").Value = -1321104.75
$ws.Range("N134").Value = -8412517.5

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap / Maple Syrup
$ws.Range("H5").Value = 20763.8
$ws.Range("I5").Value = 50252
$ws.Range("J5").Value = 1105
$ws.Range("K5").Value = 150756
$ws.Range("L5").Value = 3315
$ws.Range("M5").Value = -150644
$ws.Range("N5").Value = -3539
# Row 17: Chew the Fat / Grilled Dodo
$ws.Range("H17").Value = 1000
$ws.Range("I17").Value = 1000
$ws.Range("J17").Value = 1000
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 3000
$ws.Range("M17").Value = -2831
$ws.Range("N17").Value = -3338
# Row 18: Fisher of Men / Salt Cod
$ws.Range("H18").Value = 43332.855
$ws.Range("I18").Value = 50505
$ws.Range("K18").Value = 151515
$ws.Range("M18").Value = -151346
# Row 68: Such a Butter Face / Fermented Butter
$ws.Range("H68").Value = 1249.2195
$ws.Range("I68").Value = 939.0476
$ws.Range("J68").Value = 1356
$ws.Range("K68").Value = 2817.1428
$ws.Range("L68").Value = 4068
$ws.Range("M68").Value = -2006.1428
$ws.Range("N68").Value = -5690
# Row 71: No Margarine of Error (L) / Fermented Butter
$ws.Range("H71").Value = 1249.2195
$ws.Range("I71").Value = 939.0476
$ws.Range("J71").Value = 1356
$ws.Range("K71").Value = 8451.428400000001
$ws.Range("L71").Value = 12204
$ws.Range("M71").Value = -4395.428400000001
$ws.Range("N71").Value = -20316
# Row 113: Can't Eat Just One / Night Vinegar
$ws.Range("H113").Value = 1830.7273
$ws.Range("I113").Value = 2248.6072
$ws.Range("J113").Value = 716.381
$ws.Range("K113").Value = 6745.821599999999
$ws.Range("L113").Value = 2149.143
$ws.Range("M113").Value = -4575.821599999999
$ws.Range("N113").Value = -6489.143
# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 5684.6665
$ws.Range("I131").Value = 50489.5
$ws.Range("J131").Value = 1611.5
$ws.Range("K131").Value = 151468.5
$ws.Range("L131").Value = 4834.5
$ws.Range("M131").Value = -146428.5
$ws.Range("N131").Value = -14914.5
# Row 135: Not-so-secret Ingredient / Royal Maple Syrup
$ws.Range("H135").Value = 20763.8
$ws.Range("I135").Value = 50252
$ws.Range("J135").Value = 1105
$ws.Range("K135").Value = 452268
$ws.Range("L135").Value = 9945
$ws.Range("M135").Value = -449733
$ws.Range("N135").Value = -15015
# Row 140: Sweet, Sweet Bean Juice / Mesquite Juice
$ws.Range("H140").Value = 3300.8
$ws.Range("I140").Value = 1501.1428
$ws.Range("J140").Value = 4269.846
$ws.Range("K140").Value = 4503.428400000001
$ws.Range("L140").Value = 12809.538
$ws.Range("M140").Value = 676.5715999999993
$ws.Range("N140").Value = -23169.538

$ws = $wb.Worksheets.Item("GSM")
# Row 107: Whetstones for the Workers / Hard Mudstone Whetstone
$ws.Range("H107").Value = 4839
$ws.Range("I107").Value = 562.5
$ws.Range("J107").Value = 9115.5
$ws.Range("K107").Value = 562.5
$ws.Range("L107").Value = 9115.5
$ws.Range("M107").Value = 1357.5
$ws.Range("N107").Value = -12955.5
# Row 113: Copious Crystal Cannons / Manasilver Nugget
$ws.Range("H113").Value = 1760.3529
$ws.Range("I113").Value = 1810.2727
$ws.Range("J113").Value = 1668.8334
$ws.Range("K113").Value = 1810.2727
$ws.Range("L113").Value = 1668.8334
$ws.Range("M113").Value = 359.7273
$ws.Range("N113").Value = -6008.8334
# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 1900
$ws.Range("I122").Value = 1600
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 4800
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -2350
$ws.Range("N122").Value = -10900
# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 28575140
$ws.Range("I132").Value = 43481716
$ws.Range("J132").Value = 4202.25
$ws.Range("K132").Value = 130445148
$ws.Range("L132").Value = 12606.75
$ws.Range("M132").Value = -130442618
$ws.Range("N132").Value = -17666.75

$ws = $wb.Worksheets.Item("WVR")
# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 377211.4
$ws.Range("I136").Value = 556167.2
$ws.Range("J136").Value = 1404.25
$ws.Range("K136").Value = 1668501.6
$ws.Range("L136").Value = 4212.75
$ws.Range("M136").Value = -1665951.6
$ws.Range("N136").Value = -9312.75
